$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "45.830.42"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.89%  "
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.427.87"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.45%  "
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.37%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "115.11"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +12.08%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "319.81"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.91%  "
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.39%  "
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.633"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.35%  "
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +11.77%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0943"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.01%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.77"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.18%  "
$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.03"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.62%  "
$ws.Cells.Item(14, 2).Value = "TRON"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.110"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.40%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.05"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.70%  "
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.779.21"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.65%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.413.71"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.63%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "45.836.01"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +8.02%  "
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.46%  "
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.70%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.49"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.35%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "75.14"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.83%  "
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.56"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.15%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "269.96"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.13%  "
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.53%  "
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.64"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.56%  "
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.94%  "
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.70%  "
$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "23.13"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.82%  "
$ws.Cells.Item(31, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "39.15"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +10.09%  "
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0971"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +14.41%  "
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "172.81"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.22%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.00"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +17.45%  "
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +9.80%  "
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.29%  "
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +10.83%  "
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.13"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +13.55%  "
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +15.27%  "
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.68%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.80"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +17.11%  "
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "102.72"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -4.58%  "
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.89%  "
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "72.19"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.65%  "
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.37"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +11.67%  "
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "118.44"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +7.38%  "
$ws.Cells.Item(48, 2).Value = "THORChain"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.86"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +14.38%  "
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +17.65%  "
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.41"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +9.37%  "
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "79.10"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +3.61%  "
